# attendance.xlsx update
# - Translate the A1 header from Korean to English.
# - Add 13 new tracked columns (H:T) with their header labels.
# - Touch H:T for every row (and a handful of still-blank data cells) so
#   they serialise the same way as the rest of the sheet's "empty but
#   present" cells and the used range grows to A1:T31.
# - Fix a few existing data points (PV -> PEL / SIL typos, a stray trailing
#   quote, a literal "\n" that should be a real line break, and a doubled
#   space).
# - Replace the placeholder WO values on 2025-11-07 (row 21) with the real
#   attendance entries, and append the new rows for 2025-11-10 .. 2025-11-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. New date cells (A22:A31) need to stay plain text ("2025-11-10", ...)
#    instead of being auto-converted to date serial numbers. Clone the
#    look of the existing date column (border/bold/alignment from A2) and
#    mark the range as Text *before* typing the values.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A22:A31").PasteSpecial(-4122)
$ws.Range("A22:A31").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 1. Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Date"

# New header cells should look like the rest of row 1 (bold/bordered/
# centered style). Clone that formatting from G1 before typing the labels.
$ws.Range("G1").Copy()
$ws.Range("H1:T1").PasteSpecial(-4122)

$ws.Range("H1").Value = "사랑해"
$ws.Range("I1").Value = "ㄹㅇ.2"
$ws.Range("J1").Value = "ㄹㅇ"
$ws.Range("K1").Value = "ㄹㅇ.5"
$ws.Range("L1").Value = "ㄴㅁ.2"
$ws.Range("M1").Value = "ㄹㅇ.3"
$ws.Range("N1").Value = "ㄹㅇ.6"
$ws.Range("O1").Value = "ㄴㅁ.3"
$ws.Range("P1").Value = "ㄴㅁ.1"
$ws.Range("Q1").Value = "ㄹㅇ.4"
$ws.Range("R1").Value = "ㄹㅇ.1"
$ws.Range("S1").Value = "ㄴ"
$ws.Range("T1").Value = "ㄴㅁ"

# ---------------------------------------------------------------------------
# 2. Small corrections to existing rows (2-20)
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = "금일 오후 4시반 즈음 배터리팩 화재 발생.`n사무실내부. 초기진화 성공"
$ws.Range("F6").Value = "PEL"
$ws.Range("B7").Value = "PEL"
$ws.Range("G7").Value = "KITAS 연장으로 IMMIGRATION OFFICE 다녀옴. EGA RETURN BACK"
$ws.Range("E15").Value = "SIL"
$ws.Range("C17").Value = "SIL"
$ws.Range("G20").Value = "ㅠㅠ 설사에 몸살 기운까지.."

# ---------------------------------------------------------------------------
# 3. Row 21 (2025-11-07) - replace placeholder WO values with real entries
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "LATE(9:05)"
$ws.Range("C21").Value = "LATE(8:35)"
$ws.Range("D21").Value = "PEL"
$ws.Range("E21").Value = "LATE(8:36)"
$ws.Range("F21").Value = "ATT(8:10)"

# ---------------------------------------------------------------------------
# 4. New rows 22-31
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "2025-11-10"
$ws.Range("B22").Value = "LATE(8:38)"
$ws.Range("C22").Value = "LATE(8:50)"
$ws.Range("D22").Value = "LATE(8:37)"
$ws.Range("E22").Value = "LATE(8:35)"
$ws.Range("F22").Value = "ATT(8:12)"

$ws.Range("A23").Value = "2025-11-11"
$ws.Range("B23").Value = "LATE(9:01)"
$ws.Range("C23").Value = "WO"
$ws.Range("D23").Value = "LATE(8:36)"
$ws.Range("E23").Value = "LATE(8:33)"
$ws.Range("F23").Value = "ATT(8:15)"
$ws.Range("G23").Value = "빼빼로 day..."

$ws.Range("A24").Value = "2025-11-12"
$ws.Range("B24").Value = "WO"
$ws.Range("C24").Value = "LATE(8:40)"
$ws.Range("D24").Value = "WO"
$ws.Range("E24").Value = "LATE(8:40)"
$ws.Range("F24").Value = "ATT(8:15)"
$ws.Range("G24").Value = "수요일..."

$ws.Range("A25").Value = "2025-11-13"
$ws.Range("B25").Value = "LATE(8:33)"
$ws.Range("C25").Value = "ATT(8:27)"
$ws.Range("D25").Value = "LATE(8:36)"
$ws.Range("E25").Value = "LATE(8:33)"
$ws.Range("F25").Value = "ATT(8:15)"

$ws.Range("A26").Value = "2025-11-14"
$ws.Range("B26").Value = "LATE(8:53)"
$ws.Range("C26").Value = "LATE(8:40)"
$ws.Range("D26").Value = "LATE(8:38)"
$ws.Range("E26").Value = "LATE(8:36)"
$ws.Range("F26").Value = "ATT(8:15)"

$ws.Range("A27").Value = "2025-11-17"
$ws.Range("B27").Value = "PEL"
$ws.Range("C27").Value = "ATT(8:24)"
$ws.Range("D27").Value = "LATE(8:34)"
$ws.Range("E27").Value = "LATE(8:39)"
$ws.Range("F27").Value = "ATT(8:15)"
$ws.Range("G27").Value = "MONDAY..."

$ws.Range("A28").Value = "2025-11-18"
$ws.Range("B28").Value = "PEL"
$ws.Range("C28").Value = "ATT(8:24)"
$ws.Range("D28").Value = "LATE(8:32)"
$ws.Range("E28").Value = "LATE(8:48)"
$ws.Range("F28").Value = "ATT(8:05)"

$ws.Range("A29").Value = "2025-11-19"
$ws.Range("B29").Value = "PEL"
$ws.Range("C29").Value = "LATE(8:37)"
$ws.Range("D29").Value = "ATT(8:29)"
$ws.Range("E29").Value = "LATE(8:52)"
$ws.Range("F29").Value = "ATT(8:10)"
$ws.Range("G29").Value = "싫어욧~~ 진짜?"
$ws.Range("H29").Value = "PEL"
$ws.Range("I29").Value = "LATE(8:35)"
$ws.Range("J29").Value = "LATE(8:35)"
$ws.Range("K29").Value = "LATE(8:35)"
$ws.Range("L29").Value = "LATE(9:00)"
$ws.Range("M29").Value = "LATE(8:35)"
$ws.Range("N29").Value = "LATE(8:35)"
$ws.Range("O29").Value = "LATE(9:00)"
$ws.Range("P29").Value = "LATE(9:00)"
$ws.Range("Q29").Value = "LATE(8:35)"
$ws.Range("R29").Value = "LATE(8:35)"
$ws.Range("S29").Value = "ATT(8:21)"
$ws.Range("T29").Value = "LATE(9:00)"

$ws.Range("A30").Value = "2025-11-20"
$ws.Range("B30").Value = "PEL"
$ws.Range("C30").Value = "LATE(8:39)"
$ws.Range("D30").Value = "LATE(8:36)"
$ws.Range("E30").Value = "LATE(8:36)"
$ws.Range("F30").Value = "ATT(8:10)"

$ws.Range("A31").Value = "2025-11-21"
$ws.Range("B31").Value = "PEL"
$ws.Range("F31").Value = "ATT(8:10)"

# ---------------------------------------------------------------------------
# 5. Touch every H:T cell (rows 2-31) and the remaining blank data cells so
#    they are serialised the same way as the rest of the sheet's "empty but
#    present" cells and the used range grows to A1:T31.
# ---------------------------------------------------------------------------
$ws.Range("H2:T31").Font.Bold = $false
$ws.Range("C31:E31").Font.Bold = $false
$ws.Range("G21:G22").Font.Bold = $false
$ws.Range("G25:G26").Font.Bold = $false
$ws.Range("G28").Font.Bold = $false
$ws.Range("G30:G31").Font.Bold = $false
